$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A (was "Excel Security" for all rows 2-5) to distinct values
$ws.Range("A2").Value = "Excel Security2"
$ws.Range("A3").Value = "Excel Security3"
$ws.Range("A4").Value = "Excel Security4"
$ws.Range("A5").Value = "Excel Security5"

# Update column C values from 7 to 3
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 3
$ws.Range("C5").Value = 3

# Update the selection shown in the sheet view
$ws.Range("C2:C5").Select()
